$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set the new value in E7
$ws.Range("E7").Value = "ikljkljkljkl"

# Update the selection to match the new active cell E8
$ws.Activate()
$ws.Range("E8").Select()
